# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header style from column H and filling in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 into I1:J1
$ws.Range("H1").Copy($ws.Range("I1:J1"))

# Set the header labels
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Data for columns I (I0) and J (IF), rows 2 through 54
$iVals = @(8,1,6,6,8,7,10,7,5,6,7,3,8,7,5,5,8,7,8,6,6,1,6,7,1,1,7,1,7,7,6,6,5,7,9,5,6,8,7,8,5,6,7,7,9,7,6,5,5,6,5,4,3)
$jVals = @(8,1,6,6,9,8,10,7,6,6,7,4,8,7,5,6,8,7,8,6,6,2,6,7,1,1,7,1,7,8,6,6,6,7,9,6,7,8,8,9,5,7,8,8,9,8,6,6,6,6,5,4,3)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
